# Generate Report for Handback
#
# This applies the "handback" report update to the localization-status
# workbook:
#   - Status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it appears (Overview
#     summary sheet + the per-language detail sheets).
#   - Each per-language detail sheet (zh-cn, de-de) gets two new columns
#     filled in for every data row: F "Latest Target File" (mirrors the
#     source file name/link in column A) and G "Latest Handback File"
#     (mirrors the handoff xlf file name/link in column D).
#   - Column H "Latest Handback DateTime" is populated with the real
#     handback timestamp (was the zero-date placeholder before), and the
#     two language sheets get their own distinct timestamp.

$wb = $excel.ActiveWorkbook

$oldStatus = 'Ready for handoff'
$newStatus = 'Handed back: in sync with en-US'

# ---------------------------------------------------------------------
# 1. Flip every "Ready for handoff" status cell (Overview rollup columns
#    B/C, plus the Status column C on each language sheet) to the new
#    handed-back status text.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item('Overview')
foreach ($cellAddr in @('B2', 'C2', 'B3', 'C3')) {
    $rng = $overview.Range($cellAddr)
    if ($rng.Text -eq $oldStatus) {
        $rng.Value = $newStatus
    }
}

foreach ($sheetName in @('zh-cn', 'de-de')) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellAddr in @('C2', 'C3')) {
        $rng = $ws.Range($cellAddr)
        if ($rng.Text -eq $oldStatus) {
            $rng.Value = $newStatus
        }
    }
}

# ---------------------------------------------------------------------
# 2. Per language sheet: mirror columns A/D into new columns F/G (with
#    matching hyperlinks), and stamp the handback datetime into H.
# ---------------------------------------------------------------------
function Copy-HandbackColumn($ws, $srcAddr, $destAddr, $linkMap) {
    $display = $ws.Range($srcAddr).Text
    $ws.Range($destAddr).Value = $display
    $target = $linkMap[$srcAddr]
    if ($target) {
        $ws.Hyperlinks.Add($ws.Range($destAddr), $target, '', '', $display)
    }
}

$handbackStamps = @{ 'zh-cn' = '2016-03-17 06:05:12'; 'de-de' = '2016-03-17 06:05:27' }

foreach ($sheetName in @('zh-cn', 'de-de')) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Snapshot existing hyperlink addresses (A2/D2/A3/D3) before adding
    # any new ones, keyed by their range address.
    $linkMap = @{}
    foreach ($h in $ws.Hyperlinks) {
        $linkMap[$h.Range.Address()] = $h.Address
    }

    Copy-HandbackColumn $ws '$A$2' 'F2' $linkMap
    Copy-HandbackColumn $ws '$D$2' 'G2' $linkMap
    Copy-HandbackColumn $ws '$A$3' 'F3' $linkMap
    Copy-HandbackColumn $ws '$D$3' 'G3' $linkMap

    $stamp = $handbackStamps[$sheetName]
    $ws.Range('H2').Value = $stamp
    $ws.Range('H3').Value = $stamp
}
